$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Body content: "Version 0.0.1" change-log section
# -----------------------------------------------------------------
# Before: Initial setup for Unreal 5 (bullet) / Fixed (H2) / Changed (H2) /
#         Updated gitignore to backup correct files (bullet)
# After:  Initial setup for Unreal 5 (bullet) /
#         Basic Character Controller Motion (bullet, NEW) /
#         Fixed (H2) /
#         Debug code for line tracing to check interactions (bullet, NEW) /
#         Changed (H2) /
#         Updated gitignore to backup correct files (bullet, cleaned up)

# Locate the "Initial setup for Unreal 5" bullet paragraph - used as a
# formatting template (List Paragraph style, bulleted numId=3) for the two
# new bullet items we need to add. (NOTE: Paragraph.Range.Text includes the
# trailing paragraph mark, so trim before comparing.)
$bulletTemplate = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    $t = $pp.Range.Text.TrimEnd()
    if ($t -eq "Initial setup for Unreal 5") {
        $bulletTemplate = $pp
        break
    }
}

# --- Insert "Basic Character Controller Motion" right after that bullet ---
$srcRange = $d.Range($bulletTemplate.Range.Start, $bulletTemplate.Range.End)
$srcRange.Copy()
$insertPoint = $d.Range($bulletTemplate.Range.End, $bulletTemplate.Range.End)
$insertPoint.Paste()

$newPara1Index = $bulletTemplate.Index + 1
$newPara1 = $d.Paragraphs.Item($newPara1Index)
$r1 = $d.Range($newPara1.Range.Start, $newPara1.Range.End)
$r1.Text = "Basic Character Controller Motion"

# --- Find the "Fixed" Heading 2 paragraph (now right after the new bullet) ---
$fixedPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    $t = $pp.Range.Text.TrimEnd()
    if (($t -eq "Fixed") -and ($pp.Style.NameLocal -eq "Heading 2")) {
        $fixedPara = $pp
        break
    }
}

# --- Insert "Debug code for line tracing to check interactions" right after "Fixed" ---
$srcRange2 = $d.Range($bulletTemplate.Range.Start, $bulletTemplate.Range.End)
$srcRange2.Copy()
$insertPoint2 = $d.Range($fixedPara.Range.End, $fixedPara.Range.End)
$insertPoint2.Paste()

$newPara2Index = $fixedPara.Index + 1
$newPara2 = $d.Paragraphs.Item($newPara2Index)
$r2 = $d.Range($newPara2.Range.Start, $newPara2.Range.End)
$r2.Text = "Debug code for line tracing to check interactions"

# --- Clean up the "Updated gitignore ..." bullet: collapse the split runs
#     (with proofErr spell-check wrappers) into a single plain run ---
$gitignorePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    $t = $pp.Range.Text.TrimEnd()
    if ($t -eq "Updated gitignore to backup correct files") {
        $gitignorePara = $pp
        break
    }
}
$r3 = $d.Range($gitignorePara.Range.Start, $gitignorePara.Range.End)
$r3.Text = "Updated gitignore to backup correct files"

# -----------------------------------------------------------------
# Footer: remove the gramStart/gramEnd proofErr bookmarks around
# "HOWLONG(" and merge the "(" / "v2.0)" runs into a single "(v2.0)" run.
# -----------------------------------------------------------------
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
$ftrRange = $ftr.Range
$ftrRange.Find.Execute("HOWLONG(v2.0)", $true, $false, $false, $false, $false, $true, 1, $false, "HOWLONG(v2.0)", 2)
